$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replicate the date-column style (s="2", custom date format) used by the
# existing rows onto the newly appended date cells in column A.
$ws.Range("A357").Copy()
$ws.Range("A358:A366").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# New daily data rows (2021-08-24 .. 2021-09-01), columns A (date serial) .. AX (San Possidonio)
$data = @(
    @(44432,0,1,1,3,0,2,1,0,1,0,0,0,1,0,0,0,3,0,2,5,41,0,0,1,0,1,0,0,0,5,2,0,4,1,1,0,1,2,0,2,85,0,1,0,0,0,1,0,2),
    @(44433,1,2,0,1,1,1,0,0,0,0,0,2,0,0,0,0,2,1,0,0,22,0,2,1,2,0,0,0,0,0,0,1,1,2,0,0,0,0,0,2,44,0,0,0,0,0,0,0,0),
    @(44434,0,0,0,18,1,0,1,0,0,0,0,4,8,0,0,0,2,0,1,1,7,0,0,1,3,0,1,0,1,0,0,3,2,1,1,2,0,0,0,6,66,0,0,0,0,0,0,2,0),
    @(44435,2,1,0,11,12,2,3,2,0,0,3,2,4,0,0,0,4,1,0,3,36,0,3,1,0,0,0,1,2,5,7,9,1,0,0,1,4,8,0,5,141,1,4,0,0,1,0,1,1),
    @(44436,0,5,0,21,1,0,1,0,0,1,1,1,0,0,0,0,2,1,1,0,28,0,2,0,4,0,0,2,2,2,0,5,2,4,1,2,3,4,0,6,115,0,2,0,0,0,5,4,2),
    @(44437,7,0,5,4,13,0,1,1,0,1,10,7,3,0,0,0,0,0,0,5,23,0,3,2,8,0,0,2,2,1,1,1,5,5,0,0,1,12,0,3,129,0,1,0,0,1,0,0,1),
    @(44438,2,3,0,49,10,1,2,0,0,0,0,4,0,0,0,0,5,1,0,1,26,2,6,1,5,2,0,6,7,0,2,8,1,3,0,0,2,4,1,5,160,1,0,0,0,0,0,0,0),
    @(44439,1,0,0,9,0,2,1,2,2,3,10,10,0,0,0,0,0,0,0,4,27,0,0,0,3,1,0,0,1,4,3,2,0,0,0,0,0,5,0,8,100,0,0,0,0,0,1,0,1),
    @(44440,1,0,0,2,0,0,0,1,0,0,0,0,0,0,0,0,0,0,0,0,14,0,1,0,0,0,0,0,0,0,0,3,1,0,0,0,0,1,0,1,25,0,0,0,0,0,0,0,0)
)

$startRow = 358
$numCols = 50
for ($ri = 0; $ri -lt $data.Length; $ri++) {
    $rowVals = $data[$ri]
    $r = $startRow + $ri
    $arr = New-Object "object[,]" 1,$numCols
    for ($ci = 0; $ci -lt $numCols; $ci++) {
        $arr[0, $ci] = $rowVals[$ci]
    }
    $rng = $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, $numCols))
    $rng.Value = $arr
}

Write-Host ("Dimension now: {0}" -f $ws.UsedRange.Address())